# Updates cryptos list data (price + volume%) per the Apr 1 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.745.62"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "3.551.12"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'196.39"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'586.77"
$ws.Range("E6").Value = "  -2.94%  "

$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = "  -2.61%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.209"
$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("D10").Value = "'0.625"
$ws.Range("E10").Value = "  -3.35%  "

$ws.Range("D11").Value = "'52.64"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "'0.0000288"
$ws.Range("E12").Value = "  -5.21%  "

$ws.Range("D13").Value = "'9.24"
$ws.Range("E13").Value = "  -3.54%  "

$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").Value = "'659.90"
$ws.Range("E15").Value = "  +10.76%  "

$ws.Range("D16").Value = "69.647.11"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("D17").Value = "3.563.11"
$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").Value = "'12.55"
$ws.Range("E18").Value = "  -3.39%  "

$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").Value = "'18.43"
$ws.Range("E20").Value = "  -3.38%  "

$ws.Range("D21").Value = "'0.962"
$ws.Range("E21").Value = "  -3.59%  "

$ws.Range("D22").Value = "'18.17"
$ws.Range("E22").Value = "  +1.88%  "

$ws.Range("D23").Value = "'5.35"
$ws.Range("E23").Value = "  +3.35%  "

$ws.Range("D24").Value = "'104.74"
$ws.Range("E24").Value = "  +2.92%  "

$ws.Range("D25").Value = "'4.38"
$ws.Range("E25").Value = "  -5.44%  "

$ws.Range("D26").Value = "'2.92"
$ws.Range("E26").Value = "  -3.46%  "

$ws.Range("D27").Value = "'10.20"
$ws.Range("E27").Value = "  -5.36%  "

$ws.Range("D28").Value = "'9.59"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").Value = "'33.31"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("D30").Value = "'4.39"
$ws.Range("E30").Value = "  -5.89%  "

$ws.Range("D31").Value = "'6.80"
$ws.Range("E31").Value = "  -6.27%  "

$ws.Range("D32").Value = "'11.77"
$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("D33").Value = "'0.111"
$ws.Range("E33").Value = "  -5.72%  "

$ws.Range("D34").Value = "'61.81"
$ws.Range("E34").Value = "  -2.91%  "

$ws.Range("D35").Value = "3.754.46"
$ws.Range("E35").Value = "  -3.85%  "

$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").Value = "'3.77"
$ws.Range("E36").Value = "  +6.88%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0814"
$ws.Range("E37").Value = "  -9.01%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("D39").Value = "'511.35"
$ws.Range("E39").Value = "  -5.62%  "

$ws.Range("E40").Value = "  -6.03%  "

$ws.Range("D41").Value = "'0.371"
$ws.Range("E41").Value = "  -4.91%  "

$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("D43").Value = "'34.81"
$ws.Range("E43").Value = "  -6.03%  "

$ws.Range("D44").Value = "'0.0455"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "'3.40"
$ws.Range("E45").Value = "  -0.83%  "

$ws.Range("D46").Value = "'2.88"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("E47").Value = "  -2.57%  "

$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "'8.35"
$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("B50").Value = "Jupiter"
$ws.Range("C50").Value = "https://coinranking.com/coin/qMgTxtv34+jupiter-jup"
$ws.Range("D50").Value = "'1.77"
$ws.Range("E50").Value = "  +18.66%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000237"
$ws.Range("E51").Value = "  -5.94%  "
